$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rule "R30" row: the "From" value (C10) is corrected from 18 to 1
$ws.Range("C10").Value = 1
